$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data: extra "Grade3" column (E) and an extra student row (4) ---
$ws.Range("E1").Value = "Grade3"

$ws.Range("A4").Value = "Tymotheo"
$ws.Range("B4").Value = "M"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4

$ws.Range("E2").Value = 7
$ws.Range("E3").Value = 3
$ws.Range("E4").Value = 8

# --- Center-align the whole used block (incl. the two trailing blank columns
#     F:G down through row 7, which end up as empty, centered cells) ---
$ws.Range("A1:G7").HorizontalAlignment = -4108

# --- Header row: solid "Accent 1" theme fill ---
$ws.Range("A1:E1").Interior.ThemeColor = 5

# --- Data rows: solid "Dark 2" theme fill (lightened) ---
$ws.Range("A2:E4").Interior.ThemeColor = 3

# --- Touch F8 (no visible effect) so it materializes as a real, empty cell
#     in the sheet and the used range grows to row 8 ---
$ws.Range("F8").Font.FontStyle = "Regular"

# --- View: zoom + selection parked on F8 ---
$excel.ActiveWindow.Zoom = 89
$ws.Range("F8").Select()

# --- Print setup ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
